$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 293, pushing existing rows 293..387 down to 294..388
$ws.Rows.Item(293).Insert()

# Populate the new row 293 with its values
$ws.Cells.Item(293, 1).Value = 4
$ws.Cells.Item(293, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(293, 3).Value = 'Los Lagos'
$ws.Cells.Item(293, 4).Value = 44627
$ws.Cells.Item(293, 5).Value = 10
$ws.Cells.Item(293, 6).Value = 100112006
$ws.Cells.Item(293, 7).Value = 'Repollo'
$ws.Cells.Item(293, 8).Value = 'Crespo record'
$ws.Cells.Item(293, 9).Value = 'Primera'
$ws.Cells.Item(293, 10).Value = 300
$ws.Cells.Item(293, 11).Value = 2000
$ws.Cells.Item(293, 12).Value = 2000
$ws.Cells.Item(293, 13).Value = 2000
$ws.Cells.Item(293, 14).Value = '$/unidad'
$ws.Cells.Item(293, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(293, 16).Value = 2000
$ws.Cells.Item(293, 17).Value = 1
$ws.Cells.Item(293, 18).Value = 'Hortaliza'
